$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "face/face096.png"; $ws.Range("C2").Value = "stoppen"; $ws.Range("D2").Value = "face"
$ws.Range("B3").Value = "face/face071.png"; $ws.Range("C3").Value = "zielen"; $ws.Range("D3").Value = "face"
$ws.Range("B4").Value = "car/car079.png"; $ws.Range("C4").Value = "sparen"; $ws.Range("D4").Value = "car"
$ws.Range("B5").Value = "car/car095.png"; $ws.Range("C5").Value = "lernen"; $ws.Range("D5").Value = "car"
$ws.Range("B6").Value = "car/car069.png"; $ws.Range("C6").Value = "wehen"; $ws.Range("D6").Value = "car"
$ws.Range("B7").Value = "face/face120.png"; $ws.Range("C7").Value = "holen"; $ws.Range("D7").Value = "face"
$ws.Range("B8").Value = "face/face112.png"; $ws.Range("C8").Value = "streifen"; $ws.Range("D8").Value = "face"
$ws.Range("B9").Value = "face/face066.png"; $ws.Range("C9").Value = "wählen"; $ws.Range("D9").Value = "face"
$ws.Range("B10").Value = "face/face086.png"; $ws.Range("C10").Value = "hacken"; $ws.Range("D10").Value = "face"
$ws.Range("B11").Value = "car/car072.png"; $ws.Range("C11").Value = "mögen"; $ws.Range("D11").Value = "car"
$ws.Range("B12").Value = "face/face078.png"; $ws.Range("C12").Value = "kriegen"; $ws.Range("D12").Value = "face"
$ws.Range("B13").Value = "face/face093.png"; $ws.Range("C13").Value = "deuten"; $ws.Range("D13").Value = "face"
$ws.Range("B14").Value = "car/car109.png"; $ws.Range("C14").Value = "betteln"; $ws.Range("D14").Value = "car"
$ws.Range("B15").Value = "car/car071.png"; $ws.Range("C15").Value = "parken"; $ws.Range("D15").Value = "car"
$ws.Range("B16").Value = "face/face068.png"; $ws.Range("C16").Value = "öffnen"; $ws.Range("D16").Value = "face"
$ws.Range("B17").Value = "car/car064.png"; $ws.Range("C17").Value = "hören"; $ws.Range("D17").Value = "car"
$ws.Range("B18").Value = "face/face085.png"; $ws.Range("C18").Value = "heißen"; $ws.Range("D18").Value = "face"
$ws.Range("B19").Value = "face/face076.png"; $ws.Range("C19").Value = "bergen"; $ws.Range("D19").Value = "face"
$ws.Range("B20").Value = "car/car088.png"; $ws.Range("C20").Value = "legen"; $ws.Range("D20").Value = "car"
$ws.Range("B21").Value = "face/face101.png"; $ws.Range("C21").Value = "bauen"; $ws.Range("D21").Value = "face"
$ws.Range("B22").Value = "car/car106.png"; $ws.Range("C22").Value = "rufen"; $ws.Range("D22").Value = "car"
$ws.Range("B23").Value = "car/car115.png"; $ws.Range("C23").Value = "passen"; $ws.Range("D23").Value = "car"
$ws.Range("B24").Value = "face/face095.png"; $ws.Range("C24").Value = "zögern"; $ws.Range("D24").Value = "face"
$ws.Range("B25").Value = "face/face089.png"; $ws.Range("C25").Value = "binden"; $ws.Range("D25").Value = "face"
$ws.Range("B26").Value = "car/car108.png"; $ws.Range("C26").Value = "meinen"; $ws.Range("D26").Value = "car"
$ws.Range("B27").Value = "car/car068.png"; $ws.Range("C27").Value = "wecken"; $ws.Range("D27").Value = "car"
$ws.Range("B28").Value = "face/face075.png"; $ws.Range("C28").Value = "atmen"; $ws.Range("D28").Value = "face"
$ws.Range("B29").Value = "face/face087.png"; $ws.Range("C29").Value = "spüren"; $ws.Range("D29").Value = "face"
$ws.Range("B30").Value = "car/car092.png"; $ws.Range("C30").Value = "herrschen"; $ws.Range("D30").Value = "car"
$ws.Range("B31").Value = "car/car097.png"; $ws.Range("C31").Value = "kennen"; $ws.Range("D31").Value = "car"
$ws.Range("B32").Value = "car/car101.png"; $ws.Range("C32").Value = "kranken"; $ws.Range("D32").Value = "car"
$ws.Range("B33").Value = "car/car102.png"; $ws.Range("C33").Value = "prüfen"; $ws.Range("D33").Value = "car"
